# Update forest data - 2026-02-05 12:32
#
# The "New" sheet currently holds 3 listings (rows 2-4) that have aged out
# and need to graduate into "Previously added" (appended as rows 468-470,
# keeping their original values/links/styles). "New" is then repopulated
# with 2 freshly scraped listings.

$wb = $excel.ActiveWorkbook
$wsOld = $wb.Worksheets.Item("Previously added")
$wsNew = $wb.Worksheets.Item("New")

# A cell that already carries the plain "hyperlink text" style (small blue
# font, left/center aligned) used throughout column A - used below to wipe
# out the generic blue-underline "Hyperlink" cell style that
# Hyperlinks.Add() likes to stamp on a cell.
$linkStyleSample = $wsOld.Cells.Item(2, 1)

# A normal (non-link) data-row style sample, e.g. "9.44 ha." - plain left
# aligned 8.5pt text, used to restore a cell's look after a forced
# text-number-format round trip.
$plainStyleSample = $wsOld.Cells.Item(467, 4)

# ---------------------------------------------------------------------
# 1. Remember the hyperlink targets on "New" (row order == A2, A3, A4)
#    before anything moves/gets deleted.
# ---------------------------------------------------------------------
$moveLinks = @()
foreach ($h in $wsNew.Hyperlinks) {
    $moveLinks += $h.Address
}

# ---------------------------------------------------------------------
# 2. Copy the 3 aged-out rows straight onto the bottom of
#    "Previously added" (A468:F470) - a single-shot Range.Copy so
#    values, shared-string reuse and cell styles all come along intact.
# ---------------------------------------------------------------------
$lastRow = $wsOld.UsedRange.Rows.Count
$moveCount = $wsNew.UsedRange.Rows.Count - 1

$srcRange = $wsNew.Range("A2:F" + (1 + $moveCount))
$dstRange = $wsOld.Range("A" + ($lastRow + 1) + ":F" + ($lastRow + $moveCount))
$srcRange.Copy($dstRange)

for ($i = 0; $i -lt $moveCount; $i++) {
    $destRow = $lastRow + 1 + $i
    $aCell = $wsOld.Cells.Item($destRow, 1)
    $wsOld.Hyperlinks.Add($aCell, $moveLinks[$i])

    # Hyperlinks.Add() overwrites the cell with Excel's generic themed
    # "Hyperlink" style; paste the plain link-text format back over it so
    # column A keeps looking like the rest of the sheet.
    $linkStyleSample.Copy()
    $aCell.PasteSpecial(-4122)
}

# A freshly written "Previously added" row now has the exact target row
# styling (s=3/4/4/4/4/2) - use it as the template for the two brand new
# rows going into "New".
$rowStyleTemplate = $wsOld.Range("A" + $lastRow + ":F" + $lastRow)

# ---------------------------------------------------------------------
# 3. Clear the old "New" rows/hyperlinks, then write the two freshly
#    scraped listings into rows 2-3.
# ---------------------------------------------------------------------
$wsNew.Hyperlinks.Delete()
for ($r = 1 + $moveCount; $r -ge 2; $r--) {
    $wsNew.Rows.Item($r).Delete()
}

$newRows = @(
    @{
        A = "https://www.ss.com/msg/lv/real-estate/wood/gulbene-and-reg/gulbene/ggxmf.html"
        B = "22 000 €"
        C = "Gulbene un raj."
        D = "2 ha."
        E = "50440140001"
        F = 46058.45763888889
    },
    @{
        A = "https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/zalesjes-pag/bmhjim.html"
        B = "8 000 €"
        C = "Ludza un raj."
        D = "11 ha."
        E = "68960050098"
        F = 46058.49097222222
    }
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $destRow = 2 + $i
    $data = $newRows[$i]

    # Stamp the correct row styling (s=3/4/4/4/4/2) onto the blank row
    # before writing any values into it.
    $rowStyleTemplate.Copy()
    $wsNew.Range("A" + $destRow + ":F" + $destRow).PasteSpecial(-4122)

    $aCell = $wsNew.Cells.Item($destRow, 1)
    $aCell.Value = $data.A
    $wsNew.Cells.Item($destRow, 2).Value = $data.B
    $wsNew.Cells.Item($destRow, 3).Value = $data.C
    $wsNew.Cells.Item($destRow, 4).Value = $data.D

    # E holds a cadastre number that is all digits ("50440140001") - force
    # it to stay text (matching the scraper's original output) instead of
    # being auto-coerced into a number, then restore the plain cell style
    # via a format-only paste so no stray style gets baked onto the cell
    # permanently.
    $eCell = $wsNew.Cells.Item($destRow, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $data.E
    $plainStyleSample.Copy()
    $eCell.PasteSpecial(-4122)

    $wsNew.Cells.Item($destRow, 6).Value = $data.F

    $wsNew.Hyperlinks.Add($aCell, $data.A)
    $linkStyleSample.Copy()
    $aCell.PasteSpecial(-4122)
}

Write-Host "Moved" $moveCount "rows into 'Previously added'; wrote" $newRows.Count "new rows into 'New'."
